# MDI data roundup for 2020 accessions.xlsx
# Insert two new columns (K and L) for "Height (cm)" and "Can (cm)" = Height/Can * 100,
# ahead of the existing DBH/Dist columns (which shift right by two columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank columns at K:L -- shifts old K.. onward to M..
$ws.Columns("K:L").Insert()

# 2. Header rows --------------------------------------------------------
# Row 1 unit header -- new K1/L1 match the existing "(cm)" unit label.
$ws.Range("K1").Value = "(cm)"
$ws.Range("L1").Value = "(cm)"

# Row 2 column header -- new K2/L2 mirror Height/Can (now *100, in cm).
$ws.Range("K2").Value = "Height"
$ws.Range("L2").Value = "Can"

# 3. Body rows for groups WON (3:7), STS (15:19), SCT (21:25) -----------
# These rows get live formulas: Height(cm) = Height(m)*100, Can(cm) = Can(m)*100.
$formulaRows = @(3,4,5,6,7,15,16,17,18,19,21,22,23,24,25)
foreach ($r in $formulaRows) {
    $ws.Range("K$r").Formula = "=I$r*100"
    $ws.Range("L$r").Formula = "=J$r*100"
}

# 4. Body rows for group GOR (9:13) -- pasted/typed static values (not formulas),
#    with the "fix neighbor distance data" follow-up styled in solid black font.
$ws.Range("K9").Value = 721.1
$ws.Range("L9").Value = 792.4
$ws.Range("M9").Value = 44.577

$ws.Range("K10").Value = 746.7
$ws.Range("L10").Value = 745.2
$ws.Range("M10").Value = 49.276

$ws.Range("K11").Value = 646.1
$ws.Range("L11").Value = 509
$ws.Range("M11").Value = 33.782

$ws.Range("K12").Value = 438.3
$ws.Range("L12").Value = 679.7
$ws.Range("M12").Value = 30.734

$ws.Range("K13").Value = 685.8
$ws.Range("L13").Value = 368.8
$ws.Range("M13").Value = 34.544

# Re-style that whole 5x3 block to solid black font (matches the new font added
# for this edit), keeping the centered alignment already on the column.
$ws.Range("K9:M13").Font.Color = 0

# 5. Selection left where the author's last edit was.
$ws.Range("K21:M25").Select()
